$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: height change, and C4 becomes a text percentage instead of a numeric one
# (written before B3 so the new shared-string entries land in the same order as the target file)
$ws.Range("A4:C4").EntireRow.RowHeight = 84
$ws.Range("C4").Value = "<p align=left>65%</p>"

# Row 3: B3 text gets wrapped in <p align=left>...</p>
$ws.Range("B3").Value = "<p align=left>護盾強化<br>容易跌倒<br>MP槽不容易破裂<br>被拘束時的HP傷害<br>回避、受身的無敵時間増大<br>有毒瓦斯傷害<br>被拘束容易逃脫<br>電撃傷害<br>狀態異常容易解除度<br>火焰傷害<br>被拘束時的MP傷害<br>MP乾枯時的詠唱速度<br>狀態異常耐性</p>"

# Row 6: B6 text unchanged in content (just a shared-string reindex upstream)
$ws.Range("B6").Value = "料理全体効果%<br>詠唱中受傷時的魔力喪失量%<br>料理全體効果<br>被拘束容易逃脫"

# Update the selection to match the post-edit active cell / selection range
$ws.Range("C4:C5").Select()
